$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 64
$ws.Range("H64").Value = 3122.1191
$ws.Range("I64").Value = 2922
$ws.Range("J64").Value = 3202.1667
$ws.Range("K64").Value = 2922
$ws.Range("L64").Value = 3202.1667
$ws.Range("M64").Value = -2674
$ws.Range("N64").Value = -3698.1667

# Row 67
$ws.Range("H67").Value = 3122.1191
$ws.Range("I67").Value = 2922
$ws.Range("J67").Value = 3202.1667
$ws.Range("K67").Value = 2922
$ws.Range("L67").Value = 3202.1667
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -4918.1667

# Row 70
$ws.Range("H70").Value = 1294.2222
$ws.Range("I70").Value = 1313.64
$ws.Range("J70").Value = 1051.5
$ws.Range("K70").Value = 3940.92
$ws.Range("L70").Value = 3154.5
$ws.Range("M70").Value = -3670.92
$ws.Range("N70").Value = -3694.5

# Row 73
$ws.Range("H73").Value = 1294.2222
$ws.Range("I73").Value = 1313.64
$ws.Range("J73").Value = 1051.5
$ws.Range("K73").Value = 3940.92
$ws.Range("L73").Value = 3154.5
$ws.Range("M73").Value = -3004.92
$ws.Range("N73").Value = -5026.5

# Row 74
$ws.Range("H74").Value = 4423
$ws.Range("I74").Value = 4592.2
$ws.Range("K74").Value = 4592.2
$ws.Range("M74").Value = -3656.2

# Row 76
$ws.Range("H76").Value = 3013.4614
$ws.Range("I76").Value = 2839.5
$ws.Range("K76").Value = 2839.5
$ws.Range("M76").Value = -2524.5

# Row 77
$ws.Range("H77").Value = 4423
$ws.Range("I77").Value = 4592.2
$ws.Range("K77").Value = 22961
$ws.Range("M77").Value = -18281

# Row 79
$ws.Range("H79").Value = 3013.4614
$ws.Range("I79").Value = 2839.5
$ws.Range("K79").Value = 2839.5
$ws.Range("M79").Value = -1747.5

# Row 107
$ws.Range("H107").Value = 432.1905
$ws.Range("I107").Value = 304.66666
$ws.Range("J107").Value = 1197.3334
$ws.Range("K107").Value = 304.66666
$ws.Range("L107").Value = 1197.3334
$ws.Range("M107").Value = 1615.33334
$ws.Range("N107").Value = -5037.3334

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 13365.54
$ws.Range("I32").Value = 10420.667
$ws.Range("J32").Value = 33073.54
$ws.Range("K32").Value = 10420.667
$ws.Range("L32").Value = 33073.54
$ws.Range("M32").Value = -10133.667
$ws.Range("N32").Value = -33647.54

# Row 63
$ws.Range("H63").Value = 296674.4
$ws.Range("I63").Value = 359237.16
$ws.Range("K63").Value = 359237.16
$ws.Range("M63").Value = -358551.16

# Row 66
$ws.Range("H66").Value = 296674.4
$ws.Range("I66").Value = 359237.16
$ws.Range("K66").Value = 1796185.8
$ws.Range("M66").Value = -1792753.8

# Row 132
$ws.Range("H132").Value = 8392.865
$ws.Range("I132").Value = 8605.947
$ws.Range("J132").Value = 8167.9443
$ws.Range("K132").Value = 25817.841
$ws.Range("L132").Value = 24503.8329
$ws.Range("M132").Value = -23287.841
$ws.Range("N132").Value = -29563.8329

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 6350.0625
$ws.Range("I7").Value = 10057.1
$ws.Range("K7").Value = 10057.1
$ws.Range("M7").Value = -9944.1

# Row 94
$ws.Range("H94").Value = 1349.6154
$ws.Range("I94").Value = 1007.5
$ws.Range("J94").Value = 1642.8572
$ws.Range("K94").Value = 1007.5
$ws.Range("L94").Value = 1642.8572
$ws.Range("M94").Value = -556.5
$ws.Range("N94").Value = -2544.8572

# Row 122
$ws.Range("H122").Value = 43479350
$ws.Range("I122").Value = 55556380
$ws.Range("J122").Value = 2039.6
$ws.Range("K122").Value = 166669140
$ws.Range("L122").Value = 6118.799999999999
$ws.Range("M122").Value = -166666690
$ws.Range("N122").Value = -11018.8

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 113
$ws.Range("H113").Value = 845.3
$ws.Range("I113").Value = 631.8889
$ws.Range("J113").Value = 1019.9091
$ws.Range("K113").Value = 1895.6667
$ws.Range("L113").Value = 3059.7273
$ws.Range("M113").Value = 274.3332999999998
$ws.Range("N113").Value = -7399.7273

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 6275.657
$ws.Range("I70").Value = 7615.385
$ws.Range("K70").Value = 7615.385
$ws.Range("M70").Value = -7345.385

# Row 73
$ws.Range("H73").Value = 6275.657
$ws.Range("I73").Value = 7615.385
$ws.Range("K73").Value = 7615.385
$ws.Range("M73").Value = -6679.385

# Row 80
$ws.Range("H80").Value = 2833.6667
$ws.Range("I80").Value = 2500
$ws.Range("K80").Value = 2500
$ws.Range("M80").Value = -1502

# Row 83
$ws.Range("H83").Value = 2833.6667
$ws.Range("I83").Value = 2500
$ws.Range("K83").Value = 12500
$ws.Range("M83").Value = -7508

# Row 126
$ws.Range("H126").Value = 8930668
$ws.Range("I126").Value = 15626314
$ws.Range("J126").Value = 3140
$ws.Range("K126").Value = 46878942
$ws.Range("L126").Value = 9420
$ws.Range("M126").Value = -46876472
$ws.Range("N126").Value = -14360

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 831.6667
$ws.Range("I22").Value = 831.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 831.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -536.6667
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 831.6667
$ws.Range("I27").Value = 831.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 831.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -724.6667
$ws.Range("N27").ClearContents()

# Row 55
$ws.Range("H55").Value = 201.60976
$ws.Range("I55").Value = 220.55
$ws.Range("J55").Value = 183.57143
$ws.Range("K55").Value = 220.55
$ws.Range("L55").Value = 183.57143
$ws.Range("M55").Value = -47.55000000000001
$ws.Range("N55").Value = -529.57143

# Row 68
$ws.Range("H68").Value = 2617.7778
$ws.Range("I68").Value = 2266.6667
$ws.Range("J68").Value = 3320
$ws.Range("K68").Value = 2266.6667
$ws.Range("L68").Value = 3320
$ws.Range("M68").Value = -1517.6667
$ws.Range("N68").Value = -4818

# Row 71
$ws.Range("H71").Value = 2617.7778
$ws.Range("I71").Value = 2266.6667
$ws.Range("J71").Value = 3320
$ws.Range("K71").Value = 11333.3335
$ws.Range("L71").Value = 16600
$ws.Range("M71").Value = -7589.333500000001
$ws.Range("N71").Value = -24088

# Row 82
$ws.Range("H82").Value = 1730.2354
$ws.Range("I82").Value = 1281.1
$ws.Range("K82").Value = 1281.1
$ws.Range("M82").Value = -920.0999999999999

# Row 85
$ws.Range("H85").Value = 1730.2354
$ws.Range("I85").Value = 1281.1
$ws.Range("K85").Value = 1281.1
$ws.Range("M85").Value = -33.09999999999991

# Row 93
$ws.Range("H93").Value = 2263.25
$ws.Range("I93").Value = 1476.5
$ws.Range("K93").Value = 1476.5
$ws.Range("M93").Value = -228.5

# Row 136
$ws.Range("H136").Value = 1906.2727
$ws.Range("I136").Value = 1601.85
$ws.Range("J136").Value = 2374.6155
$ws.Range("K136").Value = 4805.549999999999
$ws.Range("L136").Value = 7123.8465
$ws.Range("M136").Value = -2255.549999999999
$ws.Range("N136").Value = -12223.8465

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 136
$ws.Range("H136").Value = 1240.561
$ws.Range("I136").Value = 1122.7632
$ws.Range("J136").Value = 2732.6667
$ws.Range("K136").Value = 3368.2896
$ws.Range("L136").Value = 8198.000100000001
$ws.Range("M136").Value = -818.2896000000001
$ws.Range("N136").Value = -13298.0001
